$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update August row label
$ws.Range("A9").Value = "August (through 08-18)"

# Update August row numeric values (B9, C9 unchanged)
$ws.Range("D9").Value = 46
$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 25
$ws.Range("G9").Value = 115
$ws.Range("H9").Value = 95

# Update Total row numeric values (B10, C10 unchanged)
$ws.Range("D10").Value = 511
$ws.Range("E10").Value = 454
$ws.Range("F10").Value = 329
$ws.Range("G10").Value = 736
$ws.Range("H10").Value = 1008
